# Natmi following Dr Hou advice
# Adds a third interacting cluster ("FAPs") to the L1cam-Ptprz1 sending/target
# cluster grid, turning the 2x (Sending x Target, 2x2 partially filled) table
# into a full 3x3 (ECs / FAPs / sCs) x (ECs / FAPs / sCs) table, and refreshes
# the associated statistics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs",  "L1cam", "Ptprz1", "ECs",  3, 1,                  19.72083766666667, 59.162513, 0.8016210077351786,  0.8016210077351787,  1, 0.3333333333333333, 0.03995766666666666, 0.119873, 0.005314930928687666, 0.005314930928687667, 0.7879986578721111, 7.091987920849,     0.004260560287097475,  0.004260560287097477),
    @("ECs",  "L1cam", "Ptprz1", "FAPs", 3, 1,                  19.72083766666667, 59.162513, 0.8016210077351786,  0.8016210077351787,  2, 0.6666666666666666, 0.03069133333333333, 0.092074, 0.004082378436578614, 0.004082378436578615, 0.6052588024402222, 5.447329221962001,  0.003272520316286511,  0.003272520316286513),
    @("ECs",  "L1cam", "Ptprz1", "sCs",  3, 1,                  19.72083766666667, 59.162513, 0.8016210077351786,  0.8016210077351787,  3, 1,                   7.447354000000001,  22.342062, 0.9906026906347337,   0.9906026906347338,   146.8680592802007,   1321.812533521806, 0.7940879271317947,    0.7940879271317949),
    @("FAPs", "L1cam", "Ptprz1", "ECs",  1, 0.3333333333333333, 0.099159,          0.297477,  0.004030657259573097, 0.004030657259573097, 1, 0.3333333333333333, 0.03995766666666666, 0.119873, 0.005314930928687666, 0.005314930928687667, 0.003962162268999999, 0.035659460421,    0.00002142266493184452, 0.00002142266493184453),
    @("FAPs", "L1cam", "Ptprz1", "FAPs", 1, 0.3333333333333333, 0.099159,          0.297477,  0.004030657259573097, 0.004030657259573097, 2, 0.6666666666666666, 0.03069133333333333, 0.092074, 0.004082378436578614, 0.004082378436578615, 0.003043321922,       0.027389897298,    0.00001645466828172026, 0.00001645466828172026),
    @("FAPs", "L1cam", "Ptprz1", "sCs",  1, 0.3333333333333333, 0.099159,          0.297477,  0.004030657259573097, 0.004030657259573097, 3, 1,                   7.447354000000001,  22.342062, 0.9906026906347337,   0.9906026906347338,   0.7384721752860001,   6.646249577574,    0.003992779926359532,   0.003992779926359533),
    @("sCs",  "L1cam", "Ptprz1", "ECs",  3, 1,                  4.781202,          14.343606, 0.1943483350052483,  0.1943483350052483,  1, 0.3333333333333333, 0.03995766666666666, 0.119873, 0.005314930928687666, 0.005314930928687667, 0.191045675782,      1.719411082038,    0.001032947976658346,   0.001032947976658346),
    @("sCs",  "L1cam", "Ptprz1", "FAPs", 3, 1,                  4.781202,          14.343606, 0.1943483350052483,  0.1943483350052483,  2, 0.6666666666666666, 0.03069133333333333, 0.092074, 0.004082378436578614, 0.004082378436578615, 0.146741464316,      1.320673178844,    0.0007934034520103821,  0.0007934034520103823),
    @("sCs",  "L1cam", "Ptprz1", "sCs",  3, 1,                  4.781202,          14.343606, 0.1943483350052483,  0.1943483350052483,  3, 1,                   7.447354000000001,  22.342062, 0.9906026906347337,   0.9906026906347338,   35.60730383950801,   320.4657345555721, 0.1925219835765795,     0.1925219835765796)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
